{"js": "// Update the date line and every \"AxB=\" multiplication prompt in the table\n// to the new values described by the commit's diff. Each old string is\n// unique in the document, so a simple search-and-replace per pair is safe.\nconst replacements = [\n  [\"2025-03-28 Friday\", \"2025-03-29 Saturday\"],\n  [\"479\u00d72=\", \"312\u00d72=\"],\n  [\"715\u00d74=\", \"356\u00d74=\"],\n  [\"481\u00d73=\", \"394\u00d76=\"],\n  [\"671\u00d72=\", \"409\u00d74=\"],\n  [\"972\u00d77=\", \"426\u00d72=\"],\n  [\"285\u00d77=\", \"490\u00d78=\"],\n  [\"403\u00d76=\", \"500\u00d79=\"],\n  [\"736\u00d75=\", \"684\u00d73=\"],\n  [\"412\u00d72=\", \"476\u00d74=\"],\n  [\"126\u00d74=\", \"300\u00d73=\"],\n  [\"139\u00d73=\", \"274\u00d77=\"],\n  [\"528\u00d76=\", \"187\u00d75=\"],\n  [\"806\u00d79=\", \"930\u00d74=\"],\n  [\"102\u00d78=\", \"905\u00d79=\"],\n  [\"721\u00d74=\", \"954\u00d76=\"],\n  [\"385\u00d74=\", \"797\u00d72=\"],\n  [\"519\u00d79=\", \"513\u00d75=\"],\n  [\"927\u00d78=\", \"400\u00d73=\"],\n  [\"320\u00d73=\", \"586\u00d78=\"],\n  [\"249\u00d72=\", \"941\u00d73=\"],\n  [\"631\u00d79=\", \"536\u00d74=\"],\n  [\"517\u00d74=\", \"392\u00d77=\"],\n  [\"550\u00d79=\", \"741\u00d73=\"],\n  [\"353\u00d74=\", \"842\u00d77=\"],\n  [\"770\u00d74=\", \"346\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"AxB=\" multiplication prompt in the table\n# to the new values described by the commit's diff. Each old string is\n# unique in the document, so Find/Replace (wdReplaceAll) per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-03-28 Friday\", \"2025-03-29 Saturday\")\n    ,@(\"479\u00d72=\", \"312\u00d72=\")\n    ,@(\"715\u00d74=\", \"356\u00d74=\")\n    ,@(\"481\u00d73=\", \"394\u00d76=\")\n    ,@(\"671\u00d72=\", \"409\u00d74=\")\n    ,@(\"972\u00d77=\", \"426\u00d72=\")\n    ,@(\"285\u00d77=\", \"490\u00d78=\")\n    ,@(\"403\u00d76=\", \"500\u00d79=\")\n    ,@(\"736\u00d75=\", \"684\u00d73=\")\n    ,@(\"412\u00d72=\", \"476\u00d74=\")\n    ,@(\"126\u00d74=\", \"300\u00d73=\")\n    ,@(\"139\u00d73=\", \"274\u00d77=\")\n    ,@(\"528\u00d76=\", \"187\u00d75=\")\n    ,@(\"806\u00d79=\", \"930\u00d74=\")\n    ,@(\"102\u00d78=\", \"905\u00d79=\")\n    ,@(\"721\u00d74=\", \"954\u00d76=\")\n    ,@(\"385\u00d74=\", \"797\u00d72=\")\n    ,@(\"519\u00d79=\", \"513\u00d75=\")\n    ,@(\"927\u00d78=\", \"400\u00d73=\")\n    ,@(\"320\u00d73=\", \"586\u00d78=\")\n    ,@(\"249\u00d72=\", \"941\u00d73=\")\n    ,@(\"631\u00d79=\", \"536\u00d74=\")\n    ,@(\"517\u00d74=\", \"392\u00d77=\")\n    ,@(\"550\u00d79=\", \"741\u00d73=\")\n    ,@(\"353\u00d74=\", \"842\u00d77=\")\n    ,@(\"770\u00d74=\", \"346\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n"}
